# Apply updated crypto price/volume figures (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.406.18"
$ws.Range("E2").Value = '  -0.13%  '

$ws.Range("D3").Value = "'1.823.25"
$ws.Range("E3").Value = '  -0.50%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = "'314.95"
$ws.Range("E5").Value = '  -0.92%  '

$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("D7").Value = "'0.5129"
$ws.Range("E7").Value = '  -3.41%  '

$ws.Range("D8").Value = "'0.3932"
$ws.Range("E8").Value = '  -3.12%  '

$ws.Range("D9").Value = "'0.07663"
$ws.Range("E9").Value = '  +1.37%  '

$ws.Range("E10").Value = '  -0.05%  '

$ws.Range("D11").Value = "'41.62"
$ws.Range("E11").Value = '  -0.93%  '

$ws.Range("D12").Value = "'21.00"
$ws.Range("E12").Value = '  +0.70%  '

$ws.Range("E13").Value = '  -1.30%  '

$ws.Range("E14").Value = '  -0.11%  '

$ws.Range("D15").Value = "'7.490"
$ws.Range("E15").Value = '  -1.68%  '

$ws.Range("D16").Value = "'1.823.93"
$ws.Range("E16").Value = '  -0.88%  '

$ws.Range("D17").Value = "'93.30"
$ws.Range("E17").Value = '  +4.09%  '

$ws.Range("D18").Value = "'0.00001096"
$ws.Range("E18").Value = '  +2.28%  '

$ws.Range("D19").Value = "'0.06660"
$ws.Range("E19").Value = '  +0.86%  '

$ws.Range("D20").Value = "'17.71"
$ws.Range("E20").Value = '  +0.80%  '

$ws.Range("E21").Value = '  -0.02%  '

$ws.Range("D22").Value = "'6.116"
$ws.Range("E22").Value = '  +0.74%  '

$ws.Range("D23").Value = "'28.414.08"
$ws.Range("E23").Value = '  -0.26%  '

$ws.Range("D24").Value = "'11.18"
$ws.Range("E24").Value = '  -1.63%  '

$ws.Range("D25").Value = "'2.256"
$ws.Range("E25").Value = '  +6.57%  '

$ws.Range("E26").Value = '  +0.71%  '

$ws.Range("D27").Value = "'156.53"
$ws.Range("E27").Value = '  -0.20%  '

$ws.Range("D28").Value = "'2.034.57"
$ws.Range("E28").Value = '  -0.75%  '

$ws.Range("D29").Value = "'2.393"
$ws.Range("E29").Value = '  -2.78%  '

$ws.Range("D30").Value = "'124.09"
$ws.Range("E30").Value = '  -0.10%  '

$ws.Range("D31").Value = "'1.108"
$ws.Range("E31").Value = '  -1.65%  '

$ws.Range("D32").Value = "'0.1094"
$ws.Range("E32").Value = '  -0.04%  '

$ws.Range("D33").Value = "'5.652"
$ws.Range("E33").Value = '  -1.01%  '

$ws.Range("D34").Value = "'3.654"
$ws.Range("E34").Value = '  -0.13%  '

$ws.Range("D35").Value = "'0.07096"
$ws.Range("E35").Value = '  -1.07%  '

$ws.Range("D36").Value = "'0.2211"
$ws.Range("E36").Value = '  -2.89%  '

$ws.Range("E37").Value = '  -0.99%  '

$ws.Range("D38").Value = "'5.170"
$ws.Range("E38").Value = '  -1.96%  '

$ws.Range("D39").Value = "'8.786"
$ws.Range("E39").Value = '  -0.33%  '

$ws.Range("D40").Value = "'0.6258"
$ws.Range("E40").Value = '  -0.41%  '

$ws.Range("E41").Value = '  -1.35%  '

$ws.Range("D42").Value = "'1.171"
$ws.Range("E42").Value = '  -2.01%  '

$ws.Range("D43").Value = "'1.0000"
$ws.Range("E43").Value = '  -0.11%  '

$ws.Range("D44").Value = "'1.391"

$ws.Range("D45").Value = "'13.39"
$ws.Range("E45").Value = '  -0.49%  '

$ws.Range("D46").Value = "'3.724"
$ws.Range("E46").Value = '  +0.29%  '

$ws.Range("D47").Value = "'0.5887"
$ws.Range("E47").Value = '  +0.51%  '

$ws.Range("E48").Value = '  -1.04%  '

$ws.Range("D49").Value = "'1.979"
$ws.Range("E49").Value = '  -0.75%  '

$ws.Range("D50").Value = "'1.196"
$ws.Range("E50").Value = '  +0.10%  '

$ws.Range("D51").Value = "'0.06900"
$ws.Range("E51").Value = '  -0.16%  '
